$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Status column updates (values reuse existing shared strings "opraveno"/"neopraveno")
$ws.Range("J2").Value = "opraveno"
$ws.Range("D3").Value = "opraveno"
$ws.Range("D5").Value = "opraveno"
$ws.Range("J5").Value = "opraveno"
$ws.Range("D6").Value = "opraveno"
$ws.Range("D11").Value = "opraveno"
$ws.Range("D12").Value = "opraveno"
$ws.Range("D13").Value = "opraveno"
$ws.Range("D14").Value = "opraveno"
$ws.Range("D15").Value = "opraveno"
$ws.Range("D17").Value = "opraveno"
$ws.Range("D18").Value = "opraveno"
$ws.Range("D19").Value = "opraveno"
$ws.Range("D20").Value = "opraveno"
$ws.Range("D21").Value = "opraveno"

# Solution (Řešení) text column - must be entered in this exact order so the
# generated shared-string table matches the original author's edit order.
$ws.Range("E3").Value = "Dle zjištěných konfliktů program vypisuje jestli CFG je nebo není daného typu."
$ws.Range("E6").Value = "funkce generování ekvivalentního nedeterministického PDA opravena"
$ws.Range("E11").Value = "přidána podpora víceznakých neterminálů"
$ws.Range("E12").Value = "Program již za žádné situace nepadá a v UI se vypisují informace, pokud se něco při parsování nepodaří."
$ws.Range("E15").Value = "Program při špatné syntaxi CFG vypíše, kterou část gramatiky se nepodařilo načíst."
$ws.Range("E17").Value = "Bug v implementaci opraven"
$ws.Range("E18").Value = "CFG je zobrazena v okně programu a uživatel již vidí s čím pracuje."
$ws.Range("E19").Value = "přidána transformace do CNF a GNF"
$ws.Range("E21").Value = "opraveno - program vypisuje parsovací tabulky zarovnané do tabulky"
$ws.Range("E20").Value = 'Program vypisuje pomocné množiny po zaškrtnutí "Show interim results"'
$ws.Range("E5").Value = "Ve složce CFG Exampels je 16 ukázkových gramatik"
$ws.Range("E14").Value = "Program zobrazuje chybovou zprávu s informacemi"

# Row 13's solution text was cleared out.
$ws.Range("E13").Value = ""

# Wrapped text in column E grew these rows, bumping their auto row height.
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(20).AutoFit()
$ws.Rows.Item(21).AutoFit()

$ws.Range("J2").Select()
